$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($row, $col, $text)
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# Apply updated coin price/volume figures (rows 2-46) and the
# row 47-51 shuffle: BabyDogeCoin inserted, RenderToken dropped.
Set-TextCell 2 4 '25.754.37'
Set-TextCell 2 5 '  -0.16%  '
Set-TextCell 3 4 '1.630.94'
Set-TextCell 3 5 '  -0.29%  '
Set-TextCell 4 5 '  -0.08%  '
Set-TextCell 5 4 '215.01'
Set-TextCell 5 5 '  -0.19%  '
Set-TextCell 6 4 '0.501'
Set-TextCell 6 5 '  -0.70%  '
Set-TextCell 7 5 '  -0.07%  '
Set-TextCell 8 4 '0.255'
Set-TextCell 8 5 '  -0.87%  '
Set-TextCell 9 5 '  -1.54%  '
Set-TextCell 10 5 '  -1.94%  '
Set-TextCell 11 5 '  +0.79%  '
Set-TextCell 12 5 '  +0.23%  '
Set-TextCell 13 4 '1.856.57'
Set-TextCell 13 5 '  -0.24%  '
Set-TextCell 14 4 '1.630.39'
Set-TextCell 14 5 '  -0.59%  '
Set-TextCell 15 4 '0.557'
Set-TextCell 15 5 '  +0.17%  '
Set-TextCell 16 4 '0.0₃0761'
Set-TextCell 16 5 '  -1.94%  '
Set-TextCell 17 4 '63.06'
Set-TextCell 17 5 '  -0.06%  '
Set-TextCell 18 4 '25.769.58'
Set-TextCell 18 5 '  -0.17%  '
Set-TextCell 19 5 '  -0.08%  '
Set-TextCell 20 5 '  -0.11%  '
Set-TextCell 21 4 '192.16'
Set-TextCell 21 5 '  -1.18%  '
Set-TextCell 22 5 '  -0.17%  '
Set-TextCell 23 5 '  +1.76%  '
Set-TextCell 24 5 '  -0.09%  '
Set-TextCell 25 5 '  +2.83%  '
Set-TextCell 26 4 '142.96'
Set-TextCell 26 5 '  +2.26%  '
Set-TextCell 27 5 '  +1.96%  '
Set-TextCell 28 5 '  +0.38%  '
Set-TextCell 29 4 '15.47'
Set-TextCell 29 5 '  -0.63%  '
Set-TextCell 30 5 '  -0.07%  '
Set-TextCell 31 4 '0.0490'
Set-TextCell 31 5 '  -0.79%  '
Set-TextCell 32 5 '  +0.14%  '
Set-TextCell 33 5 '  -0.85%  '
Set-TextCell 34 5 '  -1.40%  '
Set-TextCell 35 5 '  -0.41%  '
Set-TextCell 36 4 '0.900'
Set-TextCell 36 5 '  +0.20%  '
Set-TextCell 37 4 '1.134.47'
Set-TextCell 37 5 '  +1.94%  '
Set-TextCell 38 4 '2.52'
Set-TextCell 38 5 '  -2.14%  '
Set-TextCell 39 5 '  -1.58%  '
Set-TextCell 40 5 '  -1.26%  '
Set-TextCell 41 5 '  +0.08%  '
Set-TextCell 42 5 '  +0.86%  '
Set-TextCell 43 4 '100.51'
Set-TextCell 43 5 '  +1.22%  '
Set-TextCell 44 5 '  -0.91%  '
Set-TextCell 45 4 '0.797'
Set-TextCell 45 5 '  -0.31%  '
Set-TextCell 46 4 '1.765.98'
Set-TextCell 46 5 '  -0.11%  '
Set-TextCell 47 2 'BabyDogeCoin'
Set-TextCell 47 3 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextCell 47 4 '0.0₆0112'
Set-TextCell 47 5 '  +0.80%  '
Set-TextCell 48 2 'Aave'
Set-TextCell 48 3 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextCell 48 4 '55.29'
Set-TextCell 48 5 '  -0.54%  '
Set-TextCell 49 2 'Cronos'
Set-TextCell 49 3 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextCell 49 4 '0.0507'
Set-TextCell 49 5 '  +0.90%  '
Set-TextCell 50 2 'Mantle'
Set-TextCell 50 3 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextCell 50 4 '0.418'
Set-TextCell 50 5 '  +0.09%  '
Set-TextCell 51 2 'SynthetixNetwork'
Set-TextCell 51 3 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
Set-TextCell 51 4 '2.35'
Set-TextCell 51 5 '  -6.31%  '
